# Applies the "Team members" column addition to the REST services workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column K width (new column, ~24.43 chars as dragged by the author)
$ws.Columns.Item(11).ColumnWidth = 23.66

# Fill the K column in the same order the original author typed the values,
# so that new shared-string entries are created in the expected sequence:
#   Deeptansu, Anuj, Subham, Pratyush, Satwik, Suraj, Priyamvad, Team members
$ws.Range("K4").Value  = "Deeptansu "
$ws.Range("K5").Value  = "Deeptansu "
$ws.Range("K6").Value  = "Deeptansu "
$ws.Range("K7").Value  = "Deeptansu "

$ws.Range("K8").Value  = "Anuj"
$ws.Range("K9").Value  = "Anuj"
$ws.Range("K10").Value = "Anuj"

$ws.Range("K14").Value = "Subham"
$ws.Range("K15").Value = "Subham"
$ws.Range("K16").Value = "Subham"
$ws.Range("K17").Value = "Subham"

$ws.Range("K26").Value = "Pratyush"
$ws.Range("K27").Value = "Pratyush"
$ws.Range("K28").Value = "Pratyush"
$ws.Range("K29").Value = "Pratyush"

$ws.Range("K22").Value = "Satwik"
$ws.Range("K23").Value = "Satwik"
$ws.Range("K24").Value = "Satwik"
$ws.Range("K25").Value = "Satwik"

$ws.Range("K30").Value = "Suraj"
$ws.Range("K31").Value = "Suraj"
$ws.Range("K32").Value = "Suraj"
$ws.Range("K33").Value = "Suraj"

$ws.Range("K18").Value = "Priyamvad"
$ws.Range("K34").Value = "Priyamvad"
$ws.Range("K35").Value = "Priyamvad"
$ws.Range("K36").Value = "Priyamvad"

$ws.Range("K2").Value  = "Team members"

# Update the active selection to match the edited workbook (K4 last selected cell)
# (this also clears the old topLeftCell="A22" scroll position saved in the view)
$ws.Range("K4").Select()
